# Auto-generated Excel COM-interop script
# Updates currentAveragePrice* / LevePrice* / LeveProfit* market-data columns
# across all 8 crafting-job sheets (refreshed market board snapshot).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 117.14286
$ws.Range("I2").Value = 117.14286
$ws.Range("K2").Value = 117.14286
$ws.Range("M2").Value = -4.142859999999999
$ws.Range("H40").Value = 7673.6313
$ws.Range("I40").Value = 6021.7144
$ws.Range("K40").Value = 6021.7144
$ws.Range("M40").Value = -5846.7144
$ws.Range("H92").Value = 296.84616
$ws.Range("I92").Value = 259.9091
$ws.Range("K92").Value = 259.9091
$ws.Range("M92").Value = 988.0908999999999
$ws.Range("H116").Value = 11310.177
$ws.Range("I116").Value = 5442.636
$ws.Range("J116").Value = 22067.334
$ws.Range("K116").Value = 5442.636
$ws.Range("L116").Value = 22067.334
$ws.Range("M116").Value = -2000.636
$ws.Range("N116").Value = -28951.334
$ws.Range("H127").Value = 1504.6666
$ws.Range("I127").Value = 671.1111
$ws.Range("K127").Value = 2013.3333
$ws.Range("M127").Value = 2946.6667
$ws.Range("H131").Value = 3716.1482
$ws.Range("J131").Value = 7999.875
$ws.Range("L131").Value = 23999.625
$ws.Range("N131").Value = -34079.625
$ws.Range("H132").Value = 2856.7407
$ws.Range("I132").Value = 2744.8333
$ws.Range("K132").Value = 8234.499899999999
$ws.Range("M132").Value = -5704.499899999999
$ws.Range("H137").Value = 1914.8536
$ws.Range("J137").Value = 1981.3636
$ws.Range("L137").Value = 5944.0908
$ws.Range("N137").Value = -11044.0908
$ws.Range("H141").Value = 3551.75
$ws.Range("I141").Value = 3637.4285
$ws.Range("K141").Value = 10912.2855
$ws.Range("M141").Value = -5732.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 314.44446
$ws.Range("I4").Value = 256.33334
$ws.Range("J4").Value = 430.66666
$ws.Range("K4").Value = 256.33334
$ws.Range("L4").Value = 430.66666
$ws.Range("M4").Value = -140.33334
$ws.Range("N4").Value = -662.66666
$ws.Range("H61").Value = 3322.182
$ws.Range("I61").Value = 3322.182
$ws.Range("K61").Value = 3322.182
$ws.Range("M61").Value = -3110.182
$ws.Range("H63").Value = 5179.1763
$ws.Range("J63").Value = 9501
$ws.Range("L63").Value = 9501
$ws.Range("N63").Value = -10873
$ws.Range("H66").Value = 5179.1763
$ws.Range("J66").Value = 9501
$ws.Range("L66").Value = 47505
$ws.Range("N66").Value = -54369
$ws.Range("H110").Value = 335663.2
$ws.Range("I110").Value = 387149.78
$ws.Range("K110").Value = 387149.78
$ws.Range("M110").Value = -385104.78
$ws.Range("H132").Value = 3274
$ws.Range("I132").Value = 3362.0625
$ws.Range("K132").Value = 10086.1875
$ws.Range("M132").Value = -7556.1875
$ws.Range("H136").Value = 3322.182
$ws.Range("I136").Value = 3322.182
$ws.Range("K136").Value = 9966.545999999998
$ws.Range("M136").Value = -7416.545999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1993.6389
$ws.Range("I86").Value = 1906.3214
$ws.Range("J86").Value = 2299.25
$ws.Range("K86").Value = 1906.3214
$ws.Range("L86").Value = 2299.25
$ws.Range("M86").Value = -783.3214
$ws.Range("N86").Value = -4545.25
$ws.Range("H89").Value = 1993.6389
$ws.Range("I89").Value = 1906.3214
$ws.Range("J89").Value = 2299.25
$ws.Range("K89").Value = 9531.607
$ws.Range("L89").Value = 11496.25
$ws.Range("M89").Value = -3915.607
$ws.Range("N89").Value = -22728.25
$ws.Range("H94").Value = 934.4857
$ws.Range("I94").Value = 1023.7037
$ws.Range("K94").Value = 1023.7037
$ws.Range("M94").Value = -572.7037
$ws.Range("H105").Value = 93432.37
$ws.Range("I105").Value = 102664.5
$ws.Range("K105").Value = 102664.5
$ws.Range("M105").Value = -100917.5
$ws.Range("H134").Value = 47886.78
$ws.Range("I134").Value = 4786.476
$ws.Range("J134").Value = 500440
$ws.Range("K134").Value = 14359.428
$ws.Range("L134").Value = 1501320
$ws.Range("M134").Value = -11824.428
$ws.Range("N134").Value = -1506390

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 8109.875
$ws.Range("J14").Value = 8109.875
$ws.Range("L14").Value = 8109.875
$ws.Range("N14").Value = -8449.875
$ws.Range("H16").Value = 3417.75
$ws.Range("I16").Value = 3272.7273
$ws.Range("J16").Value = 5013
$ws.Range("K16").Value = 3272.7273
$ws.Range("L16").Value = 5013
$ws.Range("M16").Value = -2985.7273
$ws.Range("N16").Value = -5587
$ws.Range("H22").Value = 367.08334
$ws.Range("I22").Value = 382.5
$ws.Range("J22").Value = 290
$ws.Range("K22").Value = 382.5
$ws.Range("L22").Value = 290
$ws.Range("M22").Value = -32.5
$ws.Range("N22").Value = -990
$ws.Range("H26").Value = 9209.875
$ws.Range("J26").Value = 9209.875
$ws.Range("L26").Value = 9209.875
$ws.Range("N26").Value = -9783.875
$ws.Range("H70").Value = 34500
$ws.Range("J70").Value = 34500
$ws.Range("L70").Value = 34500
$ws.Range("N70").Value = -35130
$ws.Range("H73").Value = 34500
$ws.Range("J73").Value = 34500
$ws.Range("L73").Value = 34500
$ws.Range("N73").Value = -36684
$ws.Range("H113").Value = 3417.75
$ws.Range("I113").Value = 3272.7273
$ws.Range("J113").Value = 5013
$ws.Range("K113").Value = 3272.7273
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = -1102.7273
$ws.Range("N113").Value = -9353
$ws.Range("H132").Value = 1649.4242
$ws.Range("I132").Value = 1341.1538
$ws.Range("J132").Value = 2794.4285
$ws.Range("K132").Value = 4023.4614
$ws.Range("L132").Value = 8383.2855
$ws.Range("M132").Value = -1493.4614
$ws.Range("N132").Value = -13443.2855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 33333350
$ws.Range("J19").Value = 49
$ws.Range("L19").Value = 147
$ws.Range("N19").Value = -495
$ws.Range("H107").Value = 96494.82000000001
$ws.Range("J107").Value = 175920
$ws.Range("L107").Value = 527760
$ws.Range("N107").Value = -531600
$ws.Range("H132").Value = 649080.75
$ws.Range("J132").Value = 1002207.4
$ws.Range("L132").Value = 9019866.6
$ws.Range("N132").Value = -9024926.6
$ws.Range("H137").Value = 2130.238
$ws.Range("J137").Value = 1561.75
$ws.Range("L137").Value = 4685.25
$ws.Range("N137").Value = -14885.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 24011156
$ws.Range("I14").Value = 46299810
$ws.Range("K14").Value = 46299810
$ws.Range("M14").Value = -46299642
$ws.Range("H122").Value = 3777.5
$ws.Range("I122").Value = 3733.375
$ws.Range("K122").Value = 11200.125
$ws.Range("M122").Value = -8750.125
$ws.Range("H132").Value = 120558.22
$ws.Range("J132").Value = 1000000
$ws.Range("L132").Value = 3000000
$ws.Range("N132").Value = -3005060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4456.28
$ws.Range("J7").Value = 4516.6665
$ws.Range("L7").Value = 4516.6665
$ws.Range("N7").Value = -4740.6665
$ws.Range("H22").Value = 470.66666
$ws.Range("I22").Value = 470.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 470.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -175.66666
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 470.66666
$ws.Range("I27").Value = 470.66666
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 470.66666
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -363.66666
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 4375.551
$ws.Range("I40").Value = 3584.7297
$ws.Range("K40").Value = 3584.7297
$ws.Range("M40").Value = -3448.7297
$ws.Range("H42").Value = 12544.091
$ws.Range("J42").Value = 12544.091
$ws.Range("L42").Value = 12544.091
$ws.Range("N42").Value = -13670.091
$ws.Range("H46").Value = 2143.8462
$ws.Range("I46").Value = 2077.8
$ws.Range("K46").Value = 2077.8
$ws.Range("M46").Value = -1889.8
$ws.Range("H49").Value = 12544.091
$ws.Range("J49").Value = 12544.091
$ws.Range("L49").Value = 12544.091
$ws.Range("N49").Value = -12838.091
$ws.Range("H93").Value = 3218.0715
$ws.Range("I93").Value = 3054.95
$ws.Range("J93").Value = 3625.875
$ws.Range("K93").Value = 3054.95
$ws.Range("L93").Value = 3625.875
$ws.Range("M93").Value = -1806.95
$ws.Range("N93").Value = -6121.875
$ws.Range("H126").Value = 4456.28
$ws.Range("J126").Value = 4516.6665
$ws.Range("L126").Value = 13549.9995
$ws.Range("N126").Value = -18489.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3750130
$ws.Range("I4").Value = 5000106.5
$ws.Range("K4").Value = 5000106.5
$ws.Range("M4").Value = -4999993.5
$ws.Range("H70").Value = 96666.664
$ws.Range("I70").Value = 90000
$ws.Range("K70").Value = 90000
$ws.Range("M70").Value = -89685
$ws.Range("H73").Value = 96666.664
$ws.Range("I73").Value = 90000
$ws.Range("K73").Value = 90000
$ws.Range("M73").Value = -88908
$ws.Range("H75").Value = 32898
$ws.Range("I75").Value = 34490
$ws.Range("J75").Value = 32500
$ws.Range("K75").Value = 34490
$ws.Range("L75").Value = 32500
$ws.Range("M75").Value = -33554
$ws.Range("N75").Value = -34372
$ws.Range("H78").Value = 32898
$ws.Range("I78").Value = 34490
$ws.Range("J78").Value = 32500
$ws.Range("K78").Value = 103470
$ws.Range("L78").Value = 97500
$ws.Range("M78").Value = -98790
$ws.Range("N78").Value = -106860
$ws.Range("H122").Value = 29414628
$ws.Range("I122").Value = 43480384
$ws.Range("K122").Value = 130441152
$ws.Range("M122").Value = -130438702
